# Auto-generated Excel COM-interop script
# Applies the diff: recalculated ellipse-fit values across many rows in
# both the 'green' and 'blue' sheets, plus an inserted data row (frame 176)
# in the 'green' sheet (old row 44 shifts down to row 45).

$wb = $excel.ActiveWorkbook
$wsGreen = $wb.Worksheets.Item("green")
$wsBlue = $wb.Worksheets.Item("blue")

# --- Step 1: insert a new row at position 44 in 'green' sheet ---
# This shifts the former row 44 (frame 177) down to row 45, unchanged,
# and leaves a blank row 44 for the new data point (frame 176).
$wsGreen.Rows("44:44").Insert()

# --- Step 2: populate the newly inserted row 44 ---
$wsGreen.Range("A44").Value = 176
$wsGreen.Range("B44").Value = 0
$wsGreen.Range("C44").Value = 914.9709158717902
$wsGreen.Range("D44").Value = 723.6714927793219
$wsGreen.Range("E44").Value = 892.1941295863071
$wsGreen.Range("F44").Value = 722.0755153687682
$wsGreen.Range("G44").Value = 63.23577499389648
$wsGreen.Range("H44").Value = "green"

# --- Step 3: update row 43 (recalculated values; frame number 176 -> 175) ---
$wsGreen.Range("A43").Value = 175
$wsGreen.Range("B43").Value = 0
$wsGreen.Range("C43").Value = 910.834587133611
$wsGreen.Range("D43").Value = 709.8515227004976
$wsGreen.Range("E43").Value = 888.4797493632168
$wsGreen.Range("F43").Value = 709.4466226102609
$wsGreen.Range("G43").Value = 65.73482513427734
$wsGreen.Range("H43").Value = "green"

# --- Step 4: update remaining recalculated cells in 'green' sheet ---
$wsGreen.Range("C3").Value = 858.25706352359
$wsGreen.Range("D3").Value = 21.6104626205988
$wsGreen.Range("G3").Value = 54.75259017944336
$wsGreen.Range("C4").Value = 850.9871798877329
$wsGreen.Range("D4").Value = 37.48242997014471
$wsGreen.Range("G4").Value = 55.52391815185547
$wsGreen.Range("C5").Value = 843.8314662885181
$wsGreen.Range("D5").Value = 52.73091409213418
$wsGreen.Range("G5").Value = 55.68231201171875
$wsGreen.Range("C6").Value = 836.5930496921485
$wsGreen.Range("D6").Value = 68.00067473902214
$wsGreen.Range("G6").Value = 55.50335311889648
$wsGreen.Range("C7").Value = 829.0294434600368
$wsGreen.Range("D7").Value = 83.4898318373808
$wsGreen.Range("G7").Value = 55.2779655456543
$wsGreen.Range("C8").Value = 820.878029302972
$wsGreen.Range("D8").Value = 101.0311420680581
$wsGreen.Range("G8").Value = 55.71365356445312
$wsGreen.Range("C9").Value = 813.7746398587553
$wsGreen.Range("D9").Value = 116.7315511409884
$wsGreen.Range("G9").Value = 56.43114471435547
$wsGreen.Range("C10").Value = 806.5271549512943
$wsGreen.Range("D10").Value = 132.157213058623
$wsGreen.Range("G10").Value = 56.42086791992188
$wsGreen.Range("C15").Value = 775.8497839935383
$wsGreen.Range("D15").Value = 213.6443764037728
$wsGreen.Range("G15").Value = 56.8713264465332
$wsGreen.Range("C16").Value = 780.391284871287
$wsGreen.Range("D16").Value = 231.4058878985301
$wsGreen.Range("G16").Value = 56.67901992797852
$wsGreen.Range("C18").Value = 790.2939712663477
$wsGreen.Range("D18").Value = 269.2848807804018
$wsGreen.Range("G18").Value = 56.81094741821289
$wsGreen.Range("C20").Value = 799.3346240090478
$wsGreen.Range("D20").Value = 303.7301803382437
$wsGreen.Range("G20").Value = 57.33015060424805
$wsGreen.Range("D24").Value = 393.3528656252703
$wsGreen.Range("G24").Value = 58.17440795898438
$wsGreen.Range("C25").Value = 827.7015706669041
$wsGreen.Range("D25").Value = 410.0331623885846
$wsGreen.Range("G25").Value = 57.68012619018555
$wsGreen.Range("C27").Value = 838.1695755754108
$wsGreen.Range("D27").Value = 447.9413779029488
$wsGreen.Range("G27").Value = 58.58638763427734
$wsGreen.Range("C29").Value = 847.4041088092478
$wsGreen.Range("D29").Value = 482.2014118651907
$wsGreen.Range("G29").Value = 59.00302505493164
$wsGreen.Range("C32").Value = 862.7883172098053
$wsGreen.Range("D32").Value = 536.2593372609039
$wsGreen.Range("G32").Value = 60.11376953125
$wsGreen.Range("C34").Value = 871.5051405634268
$wsGreen.Range("D34").Value = 570.1917331434089
$wsGreen.Range("G34").Value = 60.44887542724609
$wsGreen.Range("C35").Value = 876.2538523496171
$wsGreen.Range("D35").Value = 585.2650798397488
$wsGreen.Range("G35").Value = 59.97185897827149
$wsGreen.Range("C38").Value = 889.7363452311699
$wsGreen.Range("D38").Value = 634.8381359604797
$wsGreen.Range("G38").Value = 64.13108825683594
$wsGreen.Range("C40").Value = 898.8032053929245
$wsGreen.Range("D40").Value = 665.5657000201273
$wsGreen.Range("G40").Value = 61.2762565612793
$wsGreen.Range("C41").Value = 902.8059579087256
$wsGreen.Range("D41").Value = 681.3479619360476
$wsGreen.Range("G41").Value = 62.93911361694336

# --- Step 5: update recalculated cells in 'blue' sheet ---
$wsBlue.Range("C2").Value = 627.7302577022728
$wsBlue.Range("D2").Value = 15.60039693697956
$wsBlue.Range("G2").Value = 55.44376373291016
$wsBlue.Range("C3").Value = 633.444971327488
$wsBlue.Range("G3").Value = 55.51166534423828
$wsBlue.Range("C4").Value = 639.328264855895
$wsBlue.Range("D4").Value = 47.80178258908358
$wsBlue.Range("G4").Value = 55.18511581420898
$wsBlue.Range("C5").Value = 645.5225561131858
$wsBlue.Range("D5").Value = 64.87613147549143
$wsBlue.Range("G5").Value = 55.22041320800781
$wsBlue.Range("C6").Value = 651.8737291523602
$wsBlue.Range("D6").Value = 82.0878029302972
$wsBlue.Range("G6").Value = 55.43023681640625
$wsBlue.Range("C7").Value = 657.6662488867408
$wsBlue.Range("D7").Value = 97.8441634826622
$wsBlue.Range("G7").Value = 55.8556022644043
$wsBlue.Range("C10").Value = 675.7412972374877
$wsBlue.Range("D10").Value = 148.1719297939857
$wsBlue.Range("G10").Value = 55.56551361083984
$wsBlue.Range("C11").Value = 682.1042590810811
$wsBlue.Range("D11").Value = 165.6552935167881
$wsBlue.Range("E11").Value = 664.1311838885583
$wsBlue.Range("G11").Value = 55.89757919311523
$wsBlue.Range("C13").Value = 686.9555788072719
$wsBlue.Range("D13").Value = 196.6480263281409
$wsBlue.Range("G13").Value = 57.03709030151367
$wsBlue.Range("C15").Value = 671.2609621180513
$wsBlue.Range("D15").Value = 227.510322819831
$wsBlue.Range("G15").Value = 56.59653854370117
$wsBlue.Range("C16").Value = 664.1033895150628
$wsBlue.Range("D16").Value = 242.1076058546433
$wsBlue.Range("G16").Value = 56.93940734863281
$wsBlue.Range("C19").Value = 642.4219191847701
$wsBlue.Range("D19").Value = 285.0608175493357
$wsBlue.Range("G19").Value = 56.95405578613281
$wsBlue.Range("C21").Value = 627.1992627463082
$wsBlue.Range("D21").Value = 315.6439914256285
$wsBlue.Range("G21").Value = 57.6081428527832
$wsBlue.Range("C23").Value = 612.5558446788177
$wsBlue.Range("D23").Value = 344.5799746044781
$wsBlue.Range("G23").Value = 57.73277282714844
$wsBlue.Range("C26").Value = 581.9907394123749
$wsBlue.Range("D26").Value = 404.0579615272651
$wsBlue.Range("G26").Value = 58.62931060791016
$wsBlue.Range("C27").Value = 574.2431824897714
$wsBlue.Range("D27").Value = 419.3535271869029
$wsBlue.Range("G27").Value = 58.50009918212891
$wsBlue.Range("C28").Value = 565.974787119626
$wsBlue.Range("D28").Value = 434.5341066619001
$wsBlue.Range("E28").Value = 543.0423886158122
$wsBlue.Range("G28").Value = 58.45481109619141
$wsBlue.Range("C29").Value = 557.5898639519544
$wsBlue.Range("D29").Value = 450.7554108593167
$wsBlue.Range("G29").Value = 58.99424362182617
$wsBlue.Range("C30").Value = 549.8415815644636
$wsBlue.Range("D30").Value = 465.7832799755325
$wsBlue.Range("G30").Value = 59.69118118286133
$wsBlue.Range("C31").Value = 542.0508594810644
$wsBlue.Range("D31").Value = 481.5829530487485
$wsBlue.Range("G31").Value = 59.87882614135742
$wsBlue.Range("C32").Value = 534.2000144951285
$wsBlue.Range("D32").Value = 495.5464291505961
$wsBlue.Range("G32").Value = 59.35115051269531
$wsBlue.Range("C34").Value = 516.6702890318707
$wsBlue.Range("D34").Value = 528.7848725268218
$wsBlue.Range("G34").Value = 60.20981216430664
$wsBlue.Range("C36").Value = 500.7255229811902
$wsBlue.Range("D36").Value = 560.5580128553228
$wsBlue.Range("G36").Value = 61.24988174438477
$wsBlue.Range("C37").Value = 491.9203742775372
$wsBlue.Range("D37").Value = 575.2321724974136
$wsBlue.Range("G37").Value = 60.45911407470703
$wsBlue.Range("C39").Value = 473.8906674822479
$wsBlue.Range("D39").Value = 609.7395899862324
$wsBlue.Range("G39").Value = 62.6220817565918
$wsBlue.Range("C40").Value = 466.1547179978416
$wsBlue.Range("D40").Value = 626.6159436353456
$wsBlue.Range("G40").Value = 63.95711517333984
$wsBlue.Range("C41").Value = 457.3392767610997
$wsBlue.Range("D41").Value = 641.5229775062668
$wsBlue.Range("G41").Value = 61.89397430419922
$wsBlue.Range("C43").Value = 438.5884125262564
$wsBlue.Range("D43").Value = 675.9574857739099
$wsBlue.Range("G43").Value = 61.79985427856445
$wsBlue.Range("C44").Value = 430.5161193834654
$wsBlue.Range("D44").Value = 692.8590946694129
$wsBlue.Range("G44").Value = 65.02391052246094
$wsBlue.Range("C45").Value = 422.2350737193472
$wsBlue.Range("D45").Value = 709.8313910498745
$wsBlue.Range("G45").Value = 65.42839813232422
$wsBlue.Range("C46").Value = 412.8737654964507
$wsBlue.Range("D46").Value = 725.4808928914147
$wsBlue.Range("G46").Value = 62.73726654052734
$wsBlue.Range("C47").Value = 403.2413601134729
$wsBlue.Range("D47").Value = 742.8314738690593
$wsBlue.Range("G47").Value = 62.04449844360352

Write-Output "done"